$d = $word.ActiveDocument

# The paragraph reads:
#   ...appeared in Court for [change of plea][ on May 28, 2022.][ ]Defendant...
# "change of plea" and " on May 28, 2022." are two separate, adjacently
# formatted runs (identical rPr). Replacing either run's text directly
# would cause this engine to coalesce it with its identically-formatted
# neighbor(s) into a single run, which would not match the target
# (the diff keeps them as two distinct <w:r> elements). To avoid that,
# temporarily bold the second run so its formatting differs from its
# neighbors while both text edits happen, then clear the bold again
# (a pure formatting change does not trigger run coalescing).

$rDate = $d.Content
$rDate.Find.Execute(" on May 28, 2022.") | Out-Null
$rDate.Bold = 1

$rReason = $d.Content
$rReason.Find.Execute("change of plea") | Out-Null
$rReason.Text = "a change of plea"

$rDate2 = $d.Content
$rDate2.Find.Execute(" on May 28, 2022.") | Out-Null
$rDate2.Text = " on May 30, 2022."

$rDate3 = $d.Content
$rDate3.Find.Execute(" on May 30, 2022.") | Out-Null
$rDate3.Bold = 0
